$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.861952666666666
$ws.Range("H2").Value = 14.585858
$ws.Range("I2").Value = 0.3995648519435639
$ws.Range("J2").Value = 0.3995648519435638
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 45.81832398232289
$ws.Range("R2").Value = 412.364915840906
$ws.Range("S2").Value = 0.02658838894581579
$ws.Range("T2").Value = 0.02658838894581579

$ws.Range("G3").Value = 4.861952666666666
$ws.Range("H3").Value = 14.585858
$ws.Range("I3").Value = 0.3995648519435639
$ws.Range("J3").Value = 0.3995648519435638
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 245.9786546946058
$ws.Range("R3").Value = 2213.807892251452
$ws.Range("S3").Value = 0.1427414967407353
$ws.Range("T3").Value = 0.1427414967407353

$ws.Range("G4").Value = 4.861952666666666
$ws.Range("H4").Value = 14.585858
$ws.Range("I4").Value = 0.3995648519435639
$ws.Range("J4").Value = 0.3995648519435638
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 130.4985176225671
$ws.Range("R4").Value = 1174.486658603104
$ws.Range("S4").Value = 0.07572833403377803
$ws.Range("T4").Value = 0.07572833403377803

$ws.Range("G5").Value = 4.861952666666666
$ws.Range("H5").Value = 14.585858
$ws.Range("I5").Value = 0.3995648519435639
$ws.Range("J5").Value = 0.3995648519435638
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 266.2528725244873
$ws.Range("R5").Value = 2396.275852720386
$ws.Range("S5").Value = 0.1545066322232348
$ws.Range("T5").Value = 0.1545066322232347

$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.04932556406896855
$ws.Range("J6").Value = 0.04932556406896854
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 5.656189887898446
$ws.Range("R6").Value = 50.90570899108601
$ws.Range("S6").Value = 0.003282288910193562
$ws.Range("T6").Value = 0.003282288910193561

$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.04932556406896855
$ws.Range("J7").Value = 0.04932556406896854
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("Q7").Value = 30.36562358455689
$ws.Range("R7").Value = 273.290612261012
$ws.Range("S7").Value = 0.01762118166434738
$ws.Range("T7").Value = 0.01762118166434738

$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.04932556406896855
$ws.Range("J8").Value = 0.04932556406896854
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 16.10980785869156
$ws.Range("R8").Value = 144.988270728224
$ws.Range("S8").Value = 0.009348526963895624
$ws.Range("T8").Value = 0.009348526963895624

$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.04932556406896855
$ws.Range("J9").Value = 0.04932556406896854
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 32.86843939944067
$ws.Range("R9").Value = 295.815954594966
$ws.Range("S9").Value = 0.01907356653053198
$ws.Range("T9").Value = 0.01907356653053198

$ws.Range("G10").Value = 4.206754333333333
$ws.Range("H10").Value = 12.620263
$ws.Range("I10").Value = 0.3457193616641432
$ws.Range("J10").Value = 0.3457193616641432
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 39.64383163994344
$ws.Range("R10").Value = 356.794484759491
$ws.Range("S10").Value = 0.02300532894550928
$ws.Range("T10").Value = 0.02300532894550928

$ws.Range("G11").Value = 4.206754333333333
$ws.Range("H11").Value = 12.620263
$ws.Range("I11").Value = 0.3457193616641432
$ws.Range("J11").Value = 0.3457193616641432
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 212.8304906459469
$ws.Range("R11").Value = 1915.474415813522
$ws.Range("S11").Value = 0.1235056059013959
$ws.Range("T11").Value = 0.1235056059013959

$ws.Range("G12").Value = 4.206754333333333
$ws.Range("H12").Value = 12.620263
$ws.Range("I12").Value = 0.3457193616641432
$ws.Range("J12").Value = 0.3457193616641432
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 112.9124946579715
$ws.Range("R12").Value = 1016.212451921744
$ws.Range("S12").Value = 0.06552315894328119
$ws.Range("T12").Value = 0.0655231589432812

$ws.Range("G13").Value = 4.206754333333333
$ws.Range("H13").Value = 12.620263
$ws.Range("I13").Value = 0.3457193616641432
$ws.Range("J13").Value = 0.3457193616641432
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 230.3725482425856
$ws.Range("R13").Value = 2073.35293418327
$ws.Range("S13").Value = 0.1336852678739569
$ws.Range("T13").Value = 0.1336852678739569

$ws.Range("G14").Value = 2.499212666666667
$ws.Range("H14").Value = 7.497638
$ws.Range("I14").Value = 0.2053902223233243
$ws.Range("J14").Value = 0.2053902223233243
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 23.55221112026289
$ws.Range("R14").Value = 211.969900082366
$ws.Range("S14").Value = 0.01366735610060981
$ws.Range("T14").Value = 0.01366735610060981

$ws.Range("G15").Value = 2.499212666666667
$ws.Range("H15").Value = 7.497638
$ws.Range("I15").Value = 0.2053902223233243
$ws.Range("J15").Value = 0.2053902223233243
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 126.4415784540858
$ws.Range("R15").Value = 1137.974206086772
$ws.Range("S15").Value = 0.07337409085843377
$ws.Range("T15").Value = 0.07337409085843377

$ws.Range("G16").Value = 2.499212666666667
$ws.Range("H16").Value = 7.497638
$ws.Range("I16").Value = 0.2053902223233243
$ws.Range("J16").Value = 0.2053902223233243
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 67.08077403952711
$ws.Range("R16").Value = 603.726966355744
$ws.Range("S16").Value = 0.03892699592498072
$ws.Range("T16").Value = 0.03892699592498073

$ws.Range("G17").Value = 2.499212666666667
$ws.Range("H17").Value = 7.497638
$ws.Range("I17").Value = 0.2053902223233243
$ws.Range("J17").Value = 0.2053902223233243
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 136.8632311276273
$ws.Range("R17").Value = 1231.769080148646
$ws.Range("S17").Value = 0.07942177943930001
$ws.Range("T17").Value = 0.07942177943930001

